$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: insert a new paragraph "//comienzo 2da etapa" right before the
# paragraph that starts with "10/04/2014 1hr. Formateo..."
# ---------------------------------------------------------------------------
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.StartsWith("10/04/2014")) {
        $targetIndex = $i
    }
}
if ($targetIndex -gt 0) {
    $p = $d.Paragraphs($targetIndex)
    $p.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs($targetIndex)
    $newPara.Range.Text = "//comienzo 2da etapa"
}

# ---------------------------------------------------------------------------
# Edit 2: "25/4/2014 1hr ." -> "25/4/2014 2hr ." (the "1" becomes "2", split
# into its own run so the surrounding proofErr wrapper still brackets it)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("1hr .", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $digit = $d.Range($start, $start + 1)
    $digit.Text = "2"
    # force a run split right after the replaced digit by round-tripping a
    # character formatting property on just that first character
    $splitPoint = $d.Range($start, $start + 1)
    $splitPoint.Font.Bold = 1
    $splitPoint.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Edit 3: append new sentences right after "... Test Case MovieClub
# formateado." (before the trailing _GoBack bookmark), split into the same
# run boundaries as the target (word boundaries around fprintf / fwritef /
# modulos).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("formateado.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insPos = $rng2.End
    $insertion = $d.Range($insPos, $insPos)
    $newText = " Crea CCD a partir de fprintf(falta testear y agregar fwritef). Crea CCD a partir de modulos."
    $insertion.Text = $newText

    # Split the newly inserted text into separate runs matching:
    #   " Crea CCD a partir de " | "fprintf" | "(" | "falta testear y agregar " |
    #   "fwritef" | "). Crea CCD a partir de " | "modulos" | "."
    $boundaries = @(22, 29, 30, 54, 61, 85, 92)
    foreach ($b in $boundaries) {
        $splitPoint2 = $d.Range($insPos, $insPos + $b)
        $splitPoint2.Font.Bold = 1
        $splitPoint2.Font.Bold = 0
    }
}
